$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 3 new rows (34-36) for reg center 10005 with user ids 110033-110035
$newRows = 34..36
$userIds = 110033, 110034, 110035

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = 10005
    $ws.Cells.Item($r, 2).Value = $userIds[$i]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Update selection/view to mimic final saved state
$ws.Range("A37:XFD1048576").Select()
